$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range('A1').Value = 'Datos actualizados a 25 de Junio de 2020 a las 11:37'

# Rebuild province/city data rows (A:E) to match the refreshed dataset.
# The city list is reordered (Soria..Lugo block moved before Gran Canaria),
# and most "Casos activos"/"Recuperados" values were reset given the new source data.
$data = @(
    @('Madrid', 68852, 19425, 68852, 40736),
    @('Cataluña', 59019, 27229, 59019, 26203),
    @('Castilla y Leon', 18907, 8267, 18907, 8716),
    @('Castilla-La Mancha', 17259, 7922, 17259, 6392),
    @('Pais Vasco', 13156, 14646, 0, 1418),
    @('Andalucia', 12679, 604, 12679, 10671),
    @('Bizkaia/Vizcaya', 10332, 0, 0, 10332),
    @('Galicia', 9041, 8409, 28, 604),
    @('Ciudad Real', 6464, 0, 0, 6464),
    @('Aragon', 5695, 1097, 5695, 3772),
    @('Valencia/Valencia', 5609, 0, 0, 5609),
    @('Zaragoza', 5287, 0, 0, 5287),
    @('Navarra', 5259, 864, 5259, 3905),
    @('Araba/Alava', 4868, 0, 0, 4868),
    @('Valladolid', 4393, 0, 0, 4393),
    @('Salamanca', 4152, 0, 0, 4152),
    @('La Rioja', 4051, 583, 4051, 3107),
    @('Toledo', 3872, 0, 0, 3872),
    @('Alacant/Alicante', 3794, 0, 0, 3794),
    @('Albacete', 3775, 0, 0, 3775),
    @('Leon', 3569, 0, 0, 3569),
    @('Segovia', 3413, 0, 0, 3413),
    @('Gipuzkoa/Guipuzcoa', 3116, 0, 0, 3116),
    @('Extremadura', 2919, 2422, 10, 487),
    @('Malaga', 2758, 0, 0, 2758),
    @('Burgos', 2746, 0, 0, 2746),
    @('Asturias', 2425, 1052, 2425, 1063),
    @('Sevilla', 2423, 0, 0, 2423),
    @('Granada', 2413, 0, 0, 2413),
    @('Soria', 2290, 0, 0, 2290),
    @('Tenerife', 2280, 0, 0, 2280),
    @('Cantabria', 2246, 1981, 62, 203),
    @('Caceres', 1973, 0, 0, 1973),
    @('A Coruña', 1969, 333, 1788, 67),
    @('Avila', 1935, 0, 0, 1935),
    @('Murcia', 1587, 2180, 0, 148),
    @('Pontevedra', 1536, 333, 1411, 30),
    @('Castello/Castellon', 1486, 0, 0, 1486),
    @('Jaen', 1387, 0, 0, 1387),
    @('Cordoba', 1331, 0, 0, 1331),
    @('Guadalajara', 1266, 0, 0, 1266),
    @('Cuenca', 1241, 0, 0, 1241),
    @('Cadiz', 1240, 0, 0, 1240),
    @('Palencia', 1205, 0, 0, 1205),
    @('Huesca', 1115, 0, 0, 1115),
    @('Zamora', 993, 0, 0, 993),
    @('Badajoz', 962, 0, 0, 962),
    @('Ourense', 751, 333, 660, 22),
    @('Teruel', 664, 0, 0, 664),
    @('Lugo', 586, 333, 520, 11),
    @('Gran Canaria', 563, 0, 0, 563),
    @('Almeria', 498, 0, 0, 498),
    @('Huelva', 400, 0, 0, 400),
    @('Mallorca', 210, 18, 194, 12),
    @('Ceuta', 125, 98, 23, 4),
    @('Melilla', 121, 125, 0, 2),
    @('La Palma', 95, 0, 0, 95),
    @('Lanzarote', 84, 0, 0, 84),
    @('Igualada, Vilanova del Cami, Santa Margarida de Montbui y Odena', 58, 0, 58, 3),
    @('Fuerteventura', 23, 0, 0, 23),
    @('Ibiza', 21, 18, 20, 1),
    @('Menorca', 15, 18, 13, 0),
    @('La Gomera', 8, 0, 0, 8),
    @('Arroyo de la Luz', 7, 0, 7, 0),
    @('El Hierro', 3, 0, 0, 3),
    @('Formentera', 0, 10, 0, 8)
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
}
